{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"Load in necessary packages\") !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (target) {\n  target.clear();\n  target.insertText(\"Load in necessary packages\", Word.InsertLocation.start);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Execute(\n    \"Load in necessary packages (tidyverse, lubridate)\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"Load in necessary packages\",\n    2\n)\n"}
